# Rewrites Sheet1's data rows (2-28) to the new matching results, and
# removes the former last row (29) entirely so the sheet ends at row 28
# (matches dimension A1:C28 instead of A1:C29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing row first so remaining writes land on the
# correct (already-shrunk) grid.
$ws.Rows.Item(29).Delete()

$ws.Cells.Item(2, 1).Value = '  main lobby cleaning'
$ws.Cells.Item(2, 2).Value = '2nd Floor Guest Lift Lobby to Clean'
$ws.Cells.Item(2, 3).Value = 0.6998999714851379

$ws.Cells.Item(3, 1).Value = ' air con grill is loose'
$ws.Cells.Item(3, 2).Value = 'Air con Grill Loose'
$ws.Cells.Item(3, 3).Value = 0.9781000018119812

$ws.Cells.Item(4, 1).Value = 'bathroom crumb polish'
$ws.Cells.Item(4, 2).Value = 'Cleaning of Bathroom'
$ws.Cells.Item(4, 3).Value = 0.6427000164985657

$ws.Cells.Item(5, 1).Value = '  fan vent in kitchen'
$ws.Cells.Item(5, 2).Value = 'Exhaust Fan'
$ws.Cells.Item(5, 3).Value = 0.6567999720573425

$ws.Cells.Item(6, 1).Value = ' blu tooth liight'
$ws.Cells.Item(6, 2).Value = 'Dental kit'
$ws.Cells.Item(6, 3).Value = 0.5746999979019165

$ws.Cells.Item(7, 1).Value = ' air band dirty'
$ws.Cells.Item(7, 2).Value = 'Air con Grill Dirty'
$ws.Cells.Item(7, 3).Value = 0.6209999918937683

$ws.Cells.Item(8, 1).Value = '  darjeeling tebags'
$ws.Cells.Item(8, 2).Value = 'Tea bags'
$ws.Cells.Item(8, 3).Value = 0.4634000062942505

$ws.Cells.Item(9, 1).Value = ' add grout to tile corner'
$ws.Cells.Item(9, 2).Value = 'Entrance Wall Socket'
$ws.Cells.Item(9, 3).Value = 0.4327999949455261

$ws.Cells.Item(10, 1).Value = '  my basin tap is loose'
$ws.Cells.Item(10, 2).Value = 'Basin Tap Loose'
$ws.Cells.Item(10, 3).Value = 0.9185000061988831

$ws.Cells.Item(11, 1).Value = 'BT'
$ws.Cells.Item(11, 2).Value = 'Express laundry service'
$ws.Cells.Item(11, 3).Value = 0.2673999965190887

$ws.Cells.Item(12, 1).Value = '  gas burner'
$ws.Cells.Item(12, 2).Value = 'Air conditioner is faulty'
$ws.Cells.Item(12, 3).Value = 0.4805000126361847

$ws.Cells.Item(13, 1).Value = ' AC duct border to repaint'
$ws.Cells.Item(13, 2).Value = 'Desk for Repaint'
$ws.Cells.Item(13, 3).Value = 0.5903000235557556

$ws.Cells.Item(14, 1).Value = '  fan body needs paint job'
$ws.Cells.Item(14, 2).Value = 'Ceiling Painting'
$ws.Cells.Item(14, 3).Value = 0.4948999881744385

$ws.Cells.Item(15, 1).Value = 'Heineken'
$ws.Cells.Item(15, 2).Value = 'Sweetener Sachet'
$ws.Cells.Item(15, 3).Value = 0.3255999982357025

$ws.Cells.Item(16, 1).Value = 'Writing Table Light Fused'
$ws.Cells.Item(16, 2).Value = 'Writing Table Light Fused'
$ws.Cells.Item(16, 3).Value = 1

$ws.Cells.Item(17, 1).Value = ' add fresh paint'
$ws.Cells.Item(17, 2).Value = 'Desk for Repaint'
$ws.Cells.Item(17, 3).Value = 0.5055999755859375

$ws.Cells.Item(18, 1).Value = 'Veuve Clicquot Brut NV  '
$ws.Cells.Item(18, 2).Value = 'Loofah'
$ws.Cells.Item(18, 3).Value = 0.2937999963760376

$ws.Cells.Item(19, 1).Value = '  hair dryer'
$ws.Cells.Item(19, 2).Value = 'Hair dryer not working'
$ws.Cells.Item(19, 3).Value = 0.9121000170707703

$ws.Cells.Item(20, 1).Value = '  pond vacuum'
$ws.Cells.Item(20, 2).Value = 'Pond Water Low Level'
$ws.Cells.Item(20, 3).Value = 0.6601999998092651

$ws.Cells.Item(21, 1).Value = 'peppermint tea'
$ws.Cells.Item(21, 2).Value = 'Chamomile tea'
$ws.Cells.Item(21, 3).Value = 0.6674000024795532

$ws.Cells.Item(22, 1).Value = ' door glass strip is broken'
$ws.Cells.Item(22, 2).Value = 'Glass Wall - Broken or Cracked'
$ws.Cells.Item(22, 3).Value = 0.7461000084877014

$ws.Cells.Item(23, 1).Value = ' bed ceiling repaint'
$ws.Cells.Item(23, 2).Value = 'Bedroom Ceiling Light Flickering'
$ws.Cells.Item(23, 3).Value = 0.6807000041007996

$ws.Cells.Item(24, 1).Value = ' buff gel'
$ws.Cells.Item(24, 2).Value = 'Bath gel'
$ws.Cells.Item(24, 3).Value = 0.6707999706268311

$ws.Cells.Item(25, 1).Value = 'set up the bait rail'
$ws.Cells.Item(25, 2).Value = 'Tape to seal box'
$ws.Cells.Item(25, 3).Value = 0.4758999943733215

$ws.Cells.Item(26, 1).Value = '  no  dnd light showing'
$ws.Cells.Item(26, 2).Value = 'Light Dimmer Not Functioning'
$ws.Cells.Item(26, 3).Value = 0.5509999990463257

$ws.Cells.Item(27, 1).Value = '  my bag rack is loose'
$ws.Cells.Item(27, 2).Value = 'Trash Bag'
$ws.Cells.Item(27, 3).Value = 0.5766000151634216

$ws.Cells.Item(28, 1).Value = '  no space bar counter'
$ws.Cells.Item(28, 2).Value = 'Bar Counter High Chair Defective'
$ws.Cells.Item(28, 3).Value = 0.6266000270843506
